$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("Tittel 1").Cut()
